# 自动更新Excel文件 - daily countdown refresh
# For each data row (2..lastRow), column E holds "days remaining" and
# column F holds the cycle start date (yyyymmdd integer), column D holds
# the total cycle length in days.
#
# Each day the sheet is refreshed:
#   - normally E is decremented by 1 (one day closer to expiry)
#   - once a cycle's last day is reached (E was 1), the cycle restarts:
#       new F = old F + D days   (the day right after the cycle ended)
#       new E = D                (fresh countdown, full length again)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # column D - total days
    $eCell = $ws.Cells.Item($r, 5)   # column E - days remaining
    $fCell = $ws.Cells.Item($r, 6)   # column F - cycle start date (yyyymmdd)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $totalDays = [int]$dVal
    $remaining = [int]$eVal
    $startSerial = [long]$fVal

    # column F is expected to be a plain yyyyMMdd integer (8 digits); rows
    # with a malformed/corrupted date (e.g. an extra stray digit) cannot be
    # parsed into a real calendar date, so leave them untouched rather than
    # crash the whole refresh.
    $year = [int]([math]::Floor($startSerial / 10000))
    $month = [int]([math]::Floor(($startSerial % 10000) / 100))
    $day = [int]($startSerial % 100)

    if ($startSerial -lt 10000101 -or $startSerial -gt 99991231 -or $month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    }
    catch {
        continue
    }

    if ($remaining -le 1) {
        # last day of the cycle -> roll over to a brand-new cycle
        $newStartDate = $startDate.AddDays($totalDays)

        $newStartSerial = [int]($newStartDate.Year * 10000 + $newStartDate.Month * 100 + $newStartDate.Day)

        $eCell.Value2 = $totalDays
        $fCell.Value2 = $newStartSerial
    }
    else {
        # one more day has passed in the current cycle
        $eCell.Value2 = $remaining - 1
    }
}
